# Updated cryptos list: refresh Price (D) and Volume(1h) (E) values for each coin row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.125.10"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "2.316.27"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "311.56"
$ws.Range("E5").Value = "  -1.41%  "
$ws.Range("D6").Value = "106.26"
$ws.Range("E6").Value = "  +3.00%  "
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "0.607"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").Value = "40.19"
$ws.Range("E10").Value = "  +1.57%  "
$ws.Range("D11").Value = "0.0914"
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("D12").Value = "8.39"
$ws.Range("E12").Value = "  -1.85%  "
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "0.994"
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").Value = "15.32"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "2.663.04"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "2.308.60"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "42.926.54"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D19").Value = "7.49"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").Value = "13.34"
$ws.Range("E21").Value = "  -3.86%  "
$ws.Range("D22").Value = "73.67"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("E23").Value = "  -1.36%  "
$ws.Range("D24").Value = "266.14"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +16.79%  "
$ws.Range("D28").Value = "10.98"
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("D30").Value = "38.79"
$ws.Range("E30").Value = "  +4.20%  "
$ws.Range("D31").Value = "22.36"
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("D32").Value = "165.63"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").Value = "0.0872"
$ws.Range("E33").Value = "  -1.28%  "
$ws.Range("D34").Value = "2.77"
$ws.Range("E34").Value = "  +8.29%  "
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("D36").Value = "4.65"
$ws.Range("E36").Value = "  +2.00%  "
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("E39").Value = "  +3.37%  "
$ws.Range("D40").Value = "3.63"
$ws.Range("E40").Value = "  -2.97%  "
$ws.Range("D41").Value = "1.59"
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("D42").Value = "104.19"
$ws.Range("E42").Value = "  +8.71%  "
$ws.Range("D43").Value = "71.02"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "12.28"
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("D47").Value = "111.93"
$ws.Range("E47").Value = "  -4.47%  "
$ws.Range("D48").Value = "1.695.48"
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("D49").Value = "76.57"
$ws.Range("E49").Value = "  -4.57%  "
$ws.Range("D50").Value = "8.88"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").Value = "5.22"
$ws.Range("E51").Value = "  -1.17%  "
